# Apply schedule updates to Sheet1 of the thesis schedule workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: add Start Date and Group assignment
$ws.Range("B5").Value = 43319
$ws.Range("F5").Value = "Group"

# Row 6: assign to Mercado, Nel
$ws.Range("F6").Value = "Mercado, Nel"

# Row 7: assign to Marzo, Pauline
$ws.Range("F7").Value = "Marzo, Pauline"

# Row 8: add End Date and Group assignment
$ws.Range("C8").Value = 43335
$ws.Range("F8").Value = "Group"

# Row 9: add End Date and Group assignment
$ws.Range("C9").Value = 43335
$ws.Range("F9").Value = "Group"

# Row 10: add End Date and assign to Alcala, Michael
$ws.Range("C10").Value = 43335
$ws.Range("F10").Value = "Alcala, Michael"

# Row 11: add End Date and assign to Mercado, Nel
$ws.Range("C11").Value = 43335
$ws.Range("F11").Value = "Mercado, Nel"

# Row 12: add End Date and assign to Garcia, Jasper
$ws.Range("C12").Value = 43335
$ws.Range("F12").Value = "Garcia, Jasper"

# Row 13: add End Date and assign to Marzo, Pauline
$ws.Range("C13").Value = 43335
$ws.Range("F13").Value = "Marzo, Pauline"

# Update the selected cell to F18 as recorded in the workbook view
$ws.Range("F18").Select()
